$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 80) mirroring the existing daily log rows.
# Column A ("日付"/date) is stored as plain text like "2025/09/22" etc. in
# the existing rows, so force the cell to text format first to stop Excel
# from auto-converting the "2025/10/08" string into a date serial number.
$rowRange = $ws.Range("A80")
$rowRange.NumberFormat = "@"
$rowRange.Value = "2025/10/08"
# Drop back to the default "Normal" style so the new cell carries no
# explicit style index, matching the unstyled data rows above it.
$rowRange.Style = "Normal"

$ws.Range("B80").Value = "水"
$ws.Range("C80").Value = 18
$ws.Range("D80").Value = 201
